# This workbook's weekly "Fruta / hortaliza" data rows (2-15) got reordered.
# Every data row's full content (A..T) moves to a different row position; the
# mapping below gives, for each destination row, which source row's data
# should end up there. We therefore snapshot all of rows 2-15 first (so we
# never read data that has already been overwritten), then write the
# snapshots back out in their new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row
$mapping = @{
    2  = 13
    3  = 5
    4  = 8
    5  = 9
    6  = 4
    7  = 10
    8  = 11
    9  = 14
    10 = 15
    11 = 12
    12 = 6
    13 = 7
    14 = 2
    15 = 3
}

$firstRow = 2
$lastRow = 15
$firstCol = 1
$lastCol = 20

# Snapshot every cell value for the affected rows before mutating anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# Now write each destination row using the snapshot of its mapped source row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcData[$c]
    }
}
